# Applies the CORP_holdings.xlsx edits:
#  - Rename fund "First Western Fixed Income Fund" -> "Oakhurst Fixed Income Fund"
#  - Update the confidential disclosure date from 2021-04-30 to 2021-05-03
#  - Refresh the Weight (D) / Percent Change (E) figures for rows 2-9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect it so the cell values can be updated.
$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect()
}

# Fund name change (row 3, column B)
$ws.Range("B3").Value = "Oakhurst Fixed Income Fund"

# Confidential disclosure text date change (row 12, column A)
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# Updated Weight / Percent Change values
$ws.Range("D2").Value = 0.1773701040128385
$ws.Range("E2").Value = 0.0009025270758122872

$ws.Range("D3").Value = 0.1772220287906256
$ws.Range("E3").Value = 0.0009823182711199419

$ws.Range("D4").Value = 0.2256366234046896
$ws.Range("E4").Value = 0.0008291873963515162

$ws.Range("D5").Value = 0.07990659254901489
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.07971249394692502
$ws.Range("E6").Value = 0

$ws.Range("D7").Value = 0.1201850540074358
$ws.Range("E7").Value = 0.0009823182711197198

$ws.Range("D8").Value = 0.1399671032884705
$ws.Range("E8").Value = 0

$ws.Range("D9").Value = 0.9999999999999998
$ws.Range("E9").Value = 0.0006393247769866939

# Restore the original protection state of the worksheet.
if ($wasProtected) {
    $ws.Protect()
}

$wb.Save()
